# updated legacy GSC export data
# The refreshed GSC export dropped the two placeholder rows for 2025-11-08 and
# 2025-11-09 (which had no real video-indexing data yet) and now starts the
# time series at 2025-11-10. Removing those two rows shifts the remaining
# data rows up by two positions and shrinks the table from 88 to 86 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the two leading placeholder rows (2025-11-08 and 2025-11-09); this
# shifts all subsequent rows up by two and re-numbers the remaining dates
# sequentially.
$ws.Rows("2:3").Delete()
